$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Statistics box updates (K/L columns, top of sheet) ---
# L7: Missing Sessions 21 -> 27
$ws.Range("L7").Value = 27
# L8: Pending Sessions 84 -> 78
$ws.Range("L8").Value = 78

# --- "Recorded By" text swap: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$gRows = @(8,9,10,12,14,15,17,34,35,36,38,40,41,43,60,61,62,64,66,67,69,86,87,88,90,92,93,95,112,113,114,116,118,119,121,138,139,140,142,144,145,147)
foreach ($r in $gRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# --- Per-group summary table (rows 21-26): P (Pending) +1, Q (Not Recorded) -1 ---
$pqRows = @(21,22,23,24,25,26)
foreach ($r in $pqRows) {
    $pCell = $ws.Range("P$r")
    $qCell = $ws.Range("Q$r")
    $pCell.Value2 = $pCell.Value2 + 1
    $qCell.Value2 = $qCell.Value2 - 1
}

# --- Rows that flip from "Pending" (yellow, style 6) to "Not Recorded" (pink, style 4) ---
$statusRows = @(177,204,231,258,285,312)
$fmtSrc = $ws.Range("A3:I3")
$fmtSrc.Copy()
foreach ($r in $statusRows) {
    $dst = $ws.Range("A$r" + ":I$r")
    $dst.PasteSpecial(-4122)
    $ws.Range("I$r").Value = "Not Recorded"
}
